$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 3, shifting existing rows 3:12 down to 5:14,
# mirroring the diff's net effect (two new weekly price records added).
$ws.Rows("3:4").Insert()

# Row 3 - new record (week of 2023-02-13)
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44970
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 100112039
$ws.Range("G3").Value = "Ciboulette"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2250
$ws.Range("N3").Value = "`$/docena de atados"
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 750
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = "Hortaliza"

# Row 4 - new record (week of 2023-02-14)
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = "Terminal La Palmera de La Serena"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44971
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 100112039
$ws.Range("G4").Value = "Ciboulette"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2250
$ws.Range("N4").Value = "`$/docena de atados"
$ws.Range("O4").Value = "Provincia del Elquí"
$ws.Range("P4").Value = 750
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = "Hortaliza"

# New D3/D4 date cells need the same number format as the other date cells
# (style index 2, numFmt "YYYY-MM-DD HH:MM:SS") - copy style from D5 which
# already carries it forward from the insert.
$ws.Range("D5").Copy()
$ws.Range("D3:D4").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
